# "fix csv and finish graben toy model"
# - remove the leftover hidden _xlchart.v1.* defined names (stale chart
#   helper ranges from an earlier chart that no longer needs them)
# - fix the header label in D1 ("sandstone" -> "formation")
# - correct/shift the B/C values for the existing fault1/fault2/limestone/
#   shale rows (2-17) and append the rest of the toy model: more shale
#   rows, the sandstone block and a brand new "shale2" block (rows 18-31)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- drop the stale hidden chart-helper defined names -----------------
$nameCount = $wb.Names.Count
for ($i = $nameCount; $i -ge 1; $i--) {
    $wb.Names.Item($i).Delete()
}

# --- header row ---------------------------------------------------------
$ws.Cells.Item(1, 4).Value = "formation"

# --- full data block (row, X, Y, Z, formation) --------------------------
$rows = @(
    @(2, 200, 430, 850, "fault1"),
    @(3, 200, 500, 500, "fault1"),
    @(4, 200, 570, 180, "fault1"),
    @(5, 200, 910, 500, "fault2"),
    @(6, 200, 980, 950, "fault2"),
    @(7, 200, 850, 180, "fault2"),
    @(8, 200, 600, 700, "limestone"),
    @(9, 200, 820, 750, "limestone"),
    @(10, 200, 1000, 920, "limestone"),
    @(11, 200, 1200, 850, "limestone"),
    @(12, 200, 400, 830, "limestone"),
    @(13, 200, 200, 770, "limestone"),
    @(14, 200, 600, 500, "shale"),
    @(15, 200, 820, 530, "shale"),
    @(16, 200, 1000, 680, "shale"),
    @(17, 200, 1200, 650, "shale"),
    @(18, 200, 400, 650, "shale"),
    @(19, 200, 200, 600, "shale"),
    @(20, 200, 200, 320, "sandstone"),
    @(21, 200, 400, 350, "sandstone"),
    @(22, 200, 600, 200, "sandstone"),
    @(23, 200, 820, 230, "sandstone"),
    @(24, 200, 1000, 370, "sandstone"),
    @(25, 200, 1200, 350, "sandstone"),
    @(26, 200, 200, 900, "shale2"),
    @(27, 200, 400, 895, "shale2"),
    @(28, 200, 600, 900, "shale2"),
    @(29, 200, 820, 905, "shale2"),
    @(30, 200, 1000, 910, "shale2"),
    @(31, 200, 1200, 910, "shale2")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# --- cosmetic: selection / view, matching the saved workbook state ------
$ws.Range("H16").Select()
